$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format Price cells that would otherwise be re-interpreted as numbers
# (losing the literal text formatting used by the source feed) as Text.
$textCells = @("D5","D8","D16","D18","D20","D21","D25","D26","D28","D30","D39","D42","D43","D44","D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.940.52"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.631.25"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "211.88"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "23.40"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "1.862.83"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "1.629.38"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "65.64"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "27.944.60"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "230.55"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "0.0₃0724"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "7.65"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  -4.94%  "
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "154.94"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("D26").Value = "6.93"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").Value = "15.56"
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("D33").Value = "1.401.83"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("E36").Value = "  +10.75%  "
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("E38").Value = "  +2.21%  "
$ws.Range("D39").Value = "0.556"
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("E40").Value = "  -3.01%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "1.85"
$ws.Range("E43").Value = "  +1.24%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "66.53"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "1.772.80"
$ws.Range("D48").Value = "88.31"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E51").Value = "  -1.21%  "
